$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Next_update (column D) for Schneider, AVM, Synology rows (2023-12-01 = 45261)
$ws.Range("D4").Value = 45261
$ws.Range("D5").Value = 45261
$ws.Range("D6").Value = 45261

# The old TP-Link row (row 7) is being pushed down to row 8 with updated values,
# and a new Swisscom row is inserted at row 7.
# First, copy current row 7 (TP-Link) data down to row 8 with updated values.
$ws.Range("A8").Value = "TP-Link"
$ws.Range("B8").Value = 100
$ws.Range("C8").Value = 44902
$ws.Range("C8").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("D8").Value = 45261
$ws.Range("D8").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("E8").Value = "TPLinkScraper"
$ws.Range("F8").Value = 20

# Now overwrite row 7 with the new Swisscom entry
$ws.Range("A7").Value = "Swisscom"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 44902
$ws.Range("D7").Value = 44927
$ws.Range("E7").Value = "SwisscomScraper"
$ws.Range("F7").ClearContents()

$ws.Range("D9").Select()
